# Generate Report for Handback
# ------------------------------------------------------------------
# This mirrors the localization tool re-running its "generate handback
# report" step: the overall status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the per-locale sheets gain links
# to the target/handback files that were produced, and the handback
# timestamps are stamped in.

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (appears on Overview!E2:F3 and on each locale sheet's Status column C2:C3)
# ------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
}

# ------------------------------------------------------------------
# 2) zh-cn sheet: fill in "Latest Target File" (I) / "Latest Handback File" (J)
#    and stamp the "Latest Handback DateTime" (K)
# ------------------------------------------------------------------
$zhcn.Range("J2").Value = "abdbf616-a93c-41f8-a0cd-5bcc5399c792.357a3f7f9a0985c4b8a6ea42f3006f4423a21e81.zh-cn.xlf"
$zhcn.Range("J3").Value = "fd65d8ed-62a2-4eed-b3db-263c3ae19e36.d3d6ac36bcc7250c1ef1c8a00655289d55b3295d.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-24 12:30:54"
$zhcn.Range("K3").Value = "2016-08-24 12:30:54"

# ------------------------------------------------------------------
# 3) de-de sheet: same, but this locale's handback ran a few seconds later
# ------------------------------------------------------------------
$dede.Range("J2").Value = "abdbf616-a93c-41f8-a0cd-5bcc5399c792.357a3f7f9a0985c4b8a6ea42f3006f4423a21e81.de-de.xlf"
$dede.Range("J3").Value = "fd65d8ed-62a2-4eed-b3db-263c3ae19e36.d3d6ac36bcc7250c1ef1c8a00655289d55b3295d.de-de.xlf"

$dede.Range("K2").Value = "2016-08-24 12:31:02"
$dede.Range("K3").Value = "2016-08-24 12:31:02"

# ------------------------------------------------------------------
# 4) Re-link "Latest Target File" (I2/I3) on both locale sheets: the
#    handback report links that column to the same source doc as
#    column A. Rebuild the hyperlink collection so the new links land
#    in doc order (A2, I2, A3, I3), matching the relationship order a
#    fresh report generation produces.
# ------------------------------------------------------------------
function Set-LocaleHyperlinks($ws) {
    $existing = @($ws.Hyperlinks)
    $linkA2 = $existing[0].Address
    $linkA3 = $existing[1].Address
    $dispA2 = $existing[0].TextToDisplay
    $dispA3 = $existing[1].TextToDisplay

    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $linkA2, [Type]::Missing, [Type]::Missing, $dispA2)
    $ws.Hyperlinks.Add($ws.Range("I2"), $linkA2, [Type]::Missing, [Type]::Missing, $dispA2)
    $ws.Hyperlinks.Add($ws.Range("A3"), $linkA3, [Type]::Missing, [Type]::Missing, $dispA3)
    $ws.Hyperlinks.Add($ws.Range("I3"), $linkA3, [Type]::Missing, [Type]::Missing, $dispA3)
}

Set-LocaleHyperlinks $zhcn
Set-LocaleHyperlinks $dede

# ------------------------------------------------------------------
# 5) Widen the columns that now hold longer text (status + the new
#    target/handback-file links) so the report stays readable.
# ------------------------------------------------------------------
$Overview.Columns.Item(5).ColumnWidth = 29.17   # E
$Overview.Columns.Item(6).ColumnWidth = 29.17   # F

$zhcn.Columns.Item(3).ColumnWidth = 29.17    # C (Status)
$zhcn.Columns.Item(9).ColumnWidth = 39.17    # I (Latest Target File)
$zhcn.Columns.Item(10).ColumnWidth = 39.17   # J (Latest Handback File)

$dede.Columns.Item(3).ColumnWidth = 29.17    # C (Status)
$dede.Columns.Item(9).ColumnWidth = 39.17    # I (Latest Target File)
$dede.Columns.Item(10).ColumnWidth = 39.17   # J (Latest Handback File)
